# Apply the crypto price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Note: Price/Volume cells are stored as plain text (no number
# formatting). Values are assigned with a leading apostrophe so Excel
# treats numeric-looking strings (e.g. "246.89") as text instead of
# silently converting them to numbers; ClearFormats() afterwards drops
# the "quote prefix" marker Excel applies to such cells, restoring the
# original default/general cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.035.24"
$ws.Range("E2").Value = "'  -2.25%  "
$ws.Range("D3").Value = "'2.241.89"
$ws.Range("E3").Value = "'  -2.15%  "
$ws.Range("E4").Value = "'  +0.26%  "
$ws.Range("D5").Value = "'246.89"
$ws.Range("E5").Value = "'  -2.03%  "
$ws.Range("D6").Value = "'0.633"
$ws.Range("E6").Value = "'  -1.46%  "
$ws.Range("D7").Value = "'77.24"
$ws.Range("E7").Value = "'  +4.46%  "
$ws.Range("E8").Value = "'  +0.17%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "'  -3.47%  "
$ws.Range("D10").Value = "'41.64"
$ws.Range("E10").Value = "'  +6.10%  "
$ws.Range("D11").Value = "'0.0956"
$ws.Range("E11").Value = "'  -2.55%  "
$ws.Range("D12").Value = "'7.14"
$ws.Range("E12").Value = "'  -4.02%  "
$ws.Range("E13").Value = "'  -2.35%  "
$ws.Range("D14").Value = "'2.581.80"
$ws.Range("E14").Value = "'  -1.89%  "
$ws.Range("D15").Value = "'14.86"
$ws.Range("E15").Value = "'  -3.41%  "
$ws.Range("D16").Value = "'0.862"
$ws.Range("E16").Value = "'  -1.64%  "
$ws.Range("D17").Value = "'2.230.95"
$ws.Range("E17").Value = "'  -2.56%  "
$ws.Range("D18").Value = "'41.955.91"
$ws.Range("E18").Value = "'  -2.12%  "
$ws.Range("D19").Value = "'0.0₃0984"
$ws.Range("E19").Value = "'  -2.36%  "
$ws.Range("D20").Value = "'6.12"
$ws.Range("E20").Value = "'  -3.03%  "
$ws.Range("D21").Value = "'71.90"
$ws.Range("E21").Value = "'  -1.11%  "
$ws.Range("D22").Value = "'2.33"
$ws.Range("E22").Value = "'  +4.06%  "
$ws.Range("D23").Value = "'231.72"
$ws.Range("E23").Value = "'  -2.72%  "
$ws.Range("E24").Value = "'  +0.04%  "
$ws.Range("E25").Value = "'  -2.24%  "
$ws.Range("E26").Value = "'  -6.47%  "
$ws.Range("E27").Value = "'  -5.26%  "
$ws.Range("D28").Value = "'7.30"
$ws.Range("E28").Value = "'  +13.35%  "
$ws.Range("E29").Value = "'  +1.20%  "
$ws.Range("D30").Value = "'170.04"
$ws.Range("E30").Value = "'  +1.74%  "
$ws.Range("D31").Value = "'20.54"
$ws.Range("E31").Value = "'  -2.55%  "
$ws.Range("E32").Value = "'  +8.86%  "
$ws.Range("D33").Value = "'0.0827"
$ws.Range("E33").Value = "'  -0.28%  "
$ws.Range("D34").Value = "'0.121"
$ws.Range("E34").Value = "'  -4.72%  "
$ws.Range("D35").Value = "'0.125"
$ws.Range("E35").Value = "'  -1.17%  "
$ws.Range("D36").Value = "'4.52"
$ws.Range("E36").Value = "'  -1.90%  "
$ws.Range("E37").Value = "'  +2.50%  "
$ws.Range("D38").Value = "'14.36"
$ws.Range("E38").Value = "'  +0.26%  "
$ws.Range("E39").Value = "'  -2.69%  "
$ws.Range("D40").Value = "'5.92"
$ws.Range("E40").Value = "'  -0.32%  "
$ws.Range("D41").Value = "'2.18"
$ws.Range("E41").Value = "'  -6.76%  "
$ws.Range("D42").Value = "'112.87"
$ws.Range("E42").Value = "'  +12.67%  "
$ws.Range("D43").Value = "'0.203"
$ws.Range("E43").Value = "'  -6.24%  "
$ws.Range("D44").Value = "'60.82"
$ws.Range("E44").Value = "'  -1.87%  "
$ws.Range("D45").Value = "'8.68"
$ws.Range("E45").Value = "'  -5.51%  "
$ws.Range("E46").Value = "'  -3.34%  "
$ws.Range("D47").Value = "'0.997"
$ws.Range("E47").Value = "'  -0.34%  "
$ws.Range("E48").Value = "'  -3.11%  "
$ws.Range("E49").Value = "'  -1.29%  "

# Rows 50 and 51 swapped rank order (FTXToken now ranks above
# NEARProtocol), along with each coin's own price/volume updates.
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'4.22"
$ws.Range("E50").Value = "'  -13.74%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'2.28"
$ws.Range("E51").Value = "'  -0.87%  "

# Remove the "quote prefix" style marker introduced by the leading
# apostrophes above so the cells' styling matches the original (no
# explicit style index / General format).
$ws.Range("D2:E51").ClearFormats()
